$wb = $excel.ActiveWorkbook

# --- Update the text summary on "Hoja1" (cell A1) with the new conversion rates ---
$wsHoja1 = $wb.Worksheets.Item("Hoja1")
$wsHoja1.Range("A1").Value = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 13.79 = 55518.9 pesos`n✅ 55518.9 pesos = 13.7 = 962.99 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"

# --- Update the rate figures on "tasas" sheet ---
$wsTasas = $wb.Worksheets.Item("tasas")
$wsTasas.Range("N10").Value = 72.5
$wsTasas.Range("O10").Value = 4025.12
$wsTasas.Range("N12").Value = 4053
$wsTasas.Range("O12").Value = 70.3
